# Updates the cryptocurrency price/volume table (columns D = Price, E = Volume(1h))
# to the latest scraped values, mirroring the automated "Updated cryptos list"
# GitHub Actions commit. Cell text is written verbatim (prices/percentages are
# display strings, not numeric types, in this sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.952.17"
# Row 3
$ws.Range("D3").Value = "1.642.99"
$ws.Range("E3").Value = "  +0.24%  "
# Row 4
$ws.Range("E4").Value = "  +0.13%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.63"
$ws.Range("E5").Value = "  +0.12%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5086"
$ws.Range("E6").Value = "  +0.95%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.23%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2562"
$ws.Range("E8").Value = "  -0.14%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06384"
$ws.Range("E9").Value = "  -0.10%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.46"
$ws.Range("E10").Value = "  -1.12%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07787"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.300"
$ws.Range("E12").Value = "  +0.94%  "
# Row 13
$ws.Range("D13").Value = "1.648.06"
$ws.Range("E13").Value = "  +0.57%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5468"
$ws.Range("E14").Value = "  +0.41%  "
# Row 15
$ws.Range("D15").Value = "0.0₅7851"
$ws.Range("E15").Value = "  -0.71%  "
# Row 16
$ws.Range("E16").Value = "  +0.17%  "
# Row 17
$ws.Range("D17").Value = "26.026.43"
$ws.Range("E17").Value = "  +0.50%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.005"
$ws.Range("E18").Value = "  +0.16%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.26"
$ws.Range("E19").Value = "  -2.47%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.448"
$ws.Range("E20").Value = "  +1.60%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.967"
$ws.Range("E21").Value = "  +0.65%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.061"
$ws.Range("E22").Value = "  +1.40%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("E23").Value = "  +0.34%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.878"
$ws.Range("E24").Value = "  -2.99%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.94"
$ws.Range("E25").Value = "  -0.06%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1148"
$ws.Range("E26").Value = "  +0.89%  "
# Row 27
$ws.Range("E27").Value = "  +1.86%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.75"
$ws.Range("E28").Value = "  +0.52%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.241"
$ws.Range("E29").Value = "  -0.22%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05030"
$ws.Range("E30").Value = "  +1.73%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.262"
$ws.Range("E31").Value = "  -0.47%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.191"
$ws.Range("E32").Value = "  +0.20%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.543"
$ws.Range("E33").Value = "  -0.10%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.369"
$ws.Range("E34").Value = "  -0.10%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8989"
$ws.Range("E35").Value = "  +0.78%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.599"
$ws.Range("E36").Value = "  -1.15%  "
# Row 37
$ws.Range("D37").Value = "1.133.32"
$ws.Range("E37").Value = "  -2.15%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5499"
$ws.Range("E38").Value = "  -1.93%  "
# Row 39
$ws.Range("E39").Value = "  +15.18%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01561"
$ws.Range("E40").Value = "  -0.28%  "
# Row 41
$ws.Range("E41").Value = "  +0.32%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.548"
$ws.Range("E42").Value = "  -0.61%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.619"
$ws.Range("E43").Value = "  -0.34%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8188"
$ws.Range("E44").Value = "  +1.55%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.26"
$ws.Range("E45").Value = "  +0.38%  "
# Row 46
$ws.Range("D46").Value = "1.779.27"
$ws.Range("E46").Value = "  +0.22%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4536"
$ws.Range("E47").Value = "  -0.07%  "
# Row 48
$ws.Range("E48").Value = "  +0.17%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.99"
$ws.Range("E49").Value = "  +0.17%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05074"
$ws.Range("E50").Value = "  +0.36%  "
# Row 51
$ws.Range("E51").Value = "  +0.59%  "
